$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the 7 rows belonging to GERARDO BAYUELO GUZMAN (73119514), which
#    was the first worker block in the table (rows 16-22). Deleting rows
#    one at a time (bottom-up) correctly shifts everything below upward,
#    which also conveniently moves the trailing "blank gap + signature"
#    rows (41-42) up to their new final position (34-35), and re-applies the
#    "closing" border style to the new last data row (29).
# ---------------------------------------------------------------------------
for ($r = 22; $r -ge 16; $r--) {
  $ws.Rows.Item($r).Delete()
}

# ---------------------------------------------------------------------------
# 2. Update the header / summary cells.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 424986      # VALOR MORA total
$ws.Range("C13").Value = 2           # Cant. Trabajadores (3 -> 2)
# F13 "Cant. Periodos" stays 7 - unchanged.

# ---------------------------------------------------------------------------
# 3. Rewrite the worker detail rows (16-29) so the two remaining workers
#    (LUIS FERNANDO MORALES CALAO and LEONARDO BAYUELO GUZMAN) are
#    interleaved period-by-period in ascending order, with the updated
#    "Salario Basico" (G) value.
# ---------------------------------------------------------------------------
$periods = @("1905", "1906", "1907", "1908", "1909", "1910", "1911")

$row = 16
foreach ($p in $periods) {
  $valorMora = 31249
  if ($p -eq "1911") { $valorMora = 24999 }

  # LUIS FERNANDO MORALES CALAO
  $ws.Range("B$row").Value = "CC"
  $ws.Range("C$row").Value = "1048441818"
  $ws.Range("D$row").Value = "LUIS FERNANDO MORALES CALAO"
  $ws.Range("E$row").Value = $p
  $ws.Range("F$row").Value = $valorMora
  $ws.Range("G$row").Value = 781241
  $row = $row + 1

  # LEONARDO BAYUELO GUZMAN
  $ws.Range("B$row").Value = "CC"
  $ws.Range("C$row").Value = "73148784"
  $ws.Range("D$row").Value = "LEONARDO BAYUELO GUZMAN"
  $ws.Range("E$row").Value = $p
  $ws.Range("F$row").Value = $valorMora
  $ws.Range("G$row").Value = 781241
  $row = $row + 1
}
